$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Rename existing sheets
# ---------------------------------------------------------------
$wb.Worksheets.Item("Spectrum").Name = "parentA"
$wb.Worksheets.Item("child1").Name = "a-1"
$wb.Worksheets.Item("child2").Name = "a-2"

# ---------------------------------------------------------------
# 2) Add the three brand-new sheets (parentB, b-1, b-2) at the end
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newParentB = $wb.Worksheets.Add($null, $lastSheet)
$newParentB.Name = "parentB"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newB1 = $wb.Worksheets.Add($null, $lastSheet)
$newB1.Name = "b-1"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newB2 = $wb.Worksheets.Add($null, $lastSheet)
$newB2.Name = "b-2"

# ---------------------------------------------------------------
# 3) Reorder the tabs into the final layout:
#    Index, Cable TV, parentA, parentB, b-1, b-2, a-1, a-2
# ---------------------------------------------------------------
$wb.Worksheets.Item("Index").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("Cable TV").Move($null, $wb.Worksheets.Item("Index"))
$wb.Worksheets.Item("parentA").Move($null, $wb.Worksheets.Item("Cable TV"))
$wb.Worksheets.Item("parentB").Move($null, $wb.Worksheets.Item("parentA"))
$wb.Worksheets.Item("b-1").Move($null, $wb.Worksheets.Item("parentB"))
$wb.Worksheets.Item("b-2").Move($null, $wb.Worksheets.Item("b-1"))
$wb.Worksheets.Item("a-1").Move($null, $wb.Worksheets.Item("b-2"))
$wb.Worksheets.Item("a-2").Move($null, $wb.Worksheets.Item("a-1"))

# ---------------------------------------------------------------
# 4) Hide the child sheets
# ---------------------------------------------------------------
$wb.Worksheets.Item("b-1").Visible = $false
$wb.Worksheets.Item("b-2").Visible = $false
$wb.Worksheets.Item("a-1").Visible = $false
$wb.Worksheets.Item("a-2").Visible = $false

# ---------------------------------------------------------------
# 5) Tab colours
#    parentA / a-1 / a-2 -> blue (FF00B0F0); parentB / b-1 / b-2 -> red (FFFF0000)
# ---------------------------------------------------------------
$wb.Worksheets.Item("parentA").Tab.Color = 15773696
$wb.Worksheets.Item("parentB").Tab.Color = 255
$wb.Worksheets.Item("b-1").Tab.Color = 255
$wb.Worksheets.Item("b-2").Tab.Color = 255

# ---------------------------------------------------------------
# 6) Cell content on the renamed / new "parent" sheets
# ---------------------------------------------------------------
$wsParentA = $wb.Worksheets.Item("parentA")
$wsParentA.Range("A1").Value = "parentA"

$wsParentB = $wb.Worksheets.Item("parentB")
$wsParentB.Range("A1").Value = "parentB"

# ---------------------------------------------------------------
# 7) Rebuild the Index sheet hyperlink list
# ---------------------------------------------------------------
$wsIndex = $wb.Worksheets.Item("Index")
$wsIndex.Hyperlinks.Delete()

$names = @("Cable TV", "parentA", "a-1", "a-2", "parentB", "b-1", "b-2")
$locations = @("'Cable TV'!A1", "parentA!A1", "'a-1'!A1", "Index!A1", "Index!A1", "Index!A1", "Index!A1")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $cell = $wsIndex.Cells.Item($row, 1)
    $cell.Value = $names[$i]
    $wsIndex.Hyperlinks.Add($cell, "", $locations[$i], "", $names[$i])
    $cell.Style = "Hyperlink"
}

# ---------------------------------------------------------------
# 8) Sheet-view selections
# ---------------------------------------------------------------
$wsParentA.Range("A2").Select()
$wsParentB.Range("A2").Select()
$wsIndex.Range("A7").Select()

$wsCableTV = $wb.Worksheets.Item("Cable TV")
$wsCableTV.Activate()
$wsCableTV.Range("B11").Select()
